# Auto-generated edit script: updates crypto price/volume cells per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.834.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.44%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.873.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.68%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.25%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'300.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.26%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.25%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.5326"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.28%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3755"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.39%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07175"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.85%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -0.02%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.8877"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.99%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08163"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.07%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.864.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.15%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'93.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.63%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'5.299"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.27%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'1.000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.34%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'14.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.44%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000008552"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.47%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -0.25%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'26.883.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.39%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -2.82%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -1.25%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'6.394"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.32%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = "'LidoDAOToken"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "'2.279"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.19%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("B25").Value = "'Monero"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "'146.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.09%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'1.734"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.52%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -1.23%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -2.82%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -2.45%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.616"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -5.70%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -1.26%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.8162"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.48%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.04985"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.88%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.177"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -4.30%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.957"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.02%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.6067"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +5.87%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.192"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -5.82%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.596"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -3.98%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -2.01%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -1.46%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'Aptos"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'8.905"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.14%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'FraxShare"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'6.578"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.43%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.5156"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +4.88%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'114.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.40%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.1494"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.58%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -0.31%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.632"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.65%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'9.905"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.84%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'37.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.70%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.06057"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.57%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'62.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -3.34%  "
$ws.Range("E51").Style = "Normal"
